$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(485).Insert()

$ws.Cells.Item(485, 1).Value = 5
$ws.Cells.Item(485, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(485, 3).Value = "Maule"
$ws.Cells.Item(485, 4).Value = 45077
$ws.Cells.Item(485, 5).Value = 7
$ws.Cells.Item(485, 6).Value = 100112023
$ws.Cells.Item(485, 7).Value = "Brócoli"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 5000
$ws.Cells.Item(485, 11).Value = 500
$ws.Cells.Item(485, 12).Value = 500
$ws.Cells.Item(485, 13).Value = 500
$ws.Cells.Item(485, 14).Value = "$/unidad"
$ws.Cells.Item(485, 15).Value = "Región del Maule"
$ws.Cells.Item(485, 16).Value = 500
$ws.Cells.Item(485, 17).Value = 1
$ws.Cells.Item(485, 18).Value = "Hortaliza"
